$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-32, matching the "repull data" commit.
$values = @{
    2  = 3
    3  = 2
    5  = -4
    6  = -1
    7  = 1
    8  = 1
    9  = -2
    10 = -3
    11 = 4
    12 = 3
    13 = 1
    14 = -1
    15 = -1
    16 = -1
    17 = 1
    18 = -1
    19 = 1
    20 = -7
    21 = -2
    24 = 2
    25 = -3
    26 = 4
    27 = 5
    29 = 5
    30 = 0
    31 = -2
    32 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
